$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 103-108 (Generation 101-106): Fitness 7293 -> 7310
for ($r = 103; $r -le 108; $r++) {
    $ws.Cells.Item($r, 3).Value = 7310
}

# Rows 109-118 (Generation 107-116): Fitness 7293 -> 7295
for ($r = 109; $r -le 118; $r++) {
    $ws.Cells.Item($r, 3).Value = 7295
}
